$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the paragraph that ends with
#    "...it just repeats the last remaining player's turn."
$d.Bookmarks.Item("_GoBack").Delete()

# 2. Re-add the "_GoBack" bookmark (collapsed) at the start of the paragraph
#    that currently reads "If you get to roll for a third time, no matter
#    what you roll, you go to jail at the end" -- do this BEFORE any of the
#    structural edits below so the bookmark tracks the paragraph correctly.
$rollParagraph = $d.Paragraphs.Item(7)
$bookmarkRange = $d.Range($rollParagraph.Range.Start, $rollParagraph.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# 3. Delete the empty paragraph that follows the "If you get to roll..."
#    paragraph, merging it away.
$emptyParagraph = $d.Paragraphs.Item(8)
$emptyParagraph.Range.Delete()

# 4. Clear out the "If you get to roll for a third time..." sentence itself,
#    leaving behind an empty paragraph that now carries the bookmark.
$rollParagraph = $d.Paragraphs.Item(7)
$textOnly = $d.Range($rollParagraph.Range.Start, $rollParagraph.Range.End - 1)
$textOnly.Delete()
